$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.333.79"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.788.41"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.20"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3945"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3452"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.10"
$ws.Range("E9").Value = "  -4.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.193"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07455"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9994"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.77"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.453"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "1.786.75"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.086"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06671"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.19"
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.68"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.509"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "27.306.60"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.381"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.477"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.517"
$ws.Range("E27").Value = "  -7.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.17"
$ws.Range("E28").Value = "  -4.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.06"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").Value = "1.989.37"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.05"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.982"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.030"
$ws.Range("E33").Value = "  -6.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08770"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.02"
$ws.Range("E35").Value = "  -6.48%  "
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.614"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.414"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6834"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02378"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06406"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2194"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.244"
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.463"
$ws.Range("E43").Value = "  -7.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.42"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9980"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6411"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.866"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.136"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.32"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07131"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.88"
$ws.Range("E51").Value = "  -2.45%  "
